$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 11d763b3...md (row 2, col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-02 13:17:17"

# zh-cn sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-02 13:17:08"
$wsZhCn.Range("K2").Value = "2016-09-02 13:17:28"

# de-de sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-02 13:17:17"
$wsDeDe.Range("K2").Value = "2016-09-02 13:17:35"
